$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated per-unit bus voltage magnitudes (vm_pu) for the 380 kV case
# Rows correspond to time steps (r=2..25); columns are bus indices (B..F, I..N).
$vmPu = @(
    @(1.02, 1.018777190134748, 1.030263256528299, 1.028441947054263, 1.035878506188936, 1.02955528997503, 1.023983701230587, 1.033074774661431, 1.031258752783233, 1.038673849366234, 1.012006349856506),
    @(1.02, 1.019837964448458, 1.030606612076188, 1.029421715039732, 1.037028821821524, 1.029546965970127, 1.024680361429499, 1.033227698754684, 1.032045993535909, 1.039632745861786, 1.012244501852445),
    @(1.02, 1.020524907074266, 1.030828240658591, 1.030056559725998, 1.037774361322411, 1.029539753632458, 1.025131167901118, 1.0333253549493, 1.032555648281803, 1.040253847108459, 1.01239841392067),
    @(1.02, 1.020813829861862, 1.030921280219724, 1.030323656259277, 1.038088075393573, 1.029536283262584, 1.025320691769214, 1.033366097555123, 1.032769968569184, 1.040515109076832, 1.012463073293099),
    @(1.02, 1.020862349024004, 1.030936894093018, 1.030368515121524, 1.038140766328743, 1.029535674837201, 1.025352513939803, 1.033372920065691, 1.032805957485931, 1.040558984963684, 1.012473927226143),
    @(1.02, 1.02052876715525, 1.030829484383945, 1.030060127866673, 1.037778552048582, 1.029539708984861, 1.025133700309516, 1.033325900582478, 1.032558511798797, 1.040257337510414, 1.012399278080471),
    @(1.02, 1.019135569909541, 1.030379406564799, 1.028772884070614, 1.036267010527895, 1.02955285428455, 1.024219136428681, 1.033126723317437, 1.031524751091472, 1.038997782203664, 1.012086873282085),
    @(1.02, 1.01668478420329, 1.029582245551909, 1.026511265982619, 1.033612736623003, 1.029562081003761, 1.022607722786081, 1.032765901988229, 1.029705111133826, 1.036783121677486, 1.011534943782861),
    @(1.02, 1.015053735864303, 1.029048220613477, 1.025008019577426, 1.03184944187669, 1.029558927936093, 1.021533564086832, 1.032518842347417, 1.028493356457054, 1.035309935747672, 1.011166039416976),
    @(1.02, 1.014348135187405, 1.028816400681362, 1.024358166769749, 1.03108739017691, 1.029555369739793, 1.021068470029094, 1.032410340936442, 1.027968971972987, 1.034672801649342, 1.011006076378505),
    @(1.02, 1.014086141422113, 1.028730207336787, 1.024116942290496, 1.0308045503853, 1.02955371971951, 1.020895716942084, 1.032369811708957, 1.027774239312679, 1.03443625652106, 1.010946625336262),
    @(1.02, 1.01414233555925, 1.028748699911505, 1.024168678527527, 1.030865210521971, 1.029554088495152, 1.020932772905002, 1.032378515607504, 1.027816007998315, 1.034486991058778, 1.010959379306702),
    @(1.02, 1.014326476710061, 1.028809277633772, 1.024338223818738, 1.03106400606181, 1.029555240040407, 1.021054190132586, 1.032406995397466, 1.027952874349895, 1.034653246413062, 1.011001162822286),
    @(1.02, 1.014439945087112, 1.028846590368057, 1.024442707445325, 1.03118651971999, 1.029555906069241, 1.021128999764434, 1.032424512714591, 1.028037208466125, 1.034755697038808, 1.011026902571408),
    @(1.02, 1.015100577992576, 1.029063593663635, 1.025051170569725, 1.031900047620276, 1.029559117975007, 1.02156443134678, 1.032526011298701, 1.028528164730228, 1.035352236350497, 1.011176650909739),
    @(1.02, 1.01551515029314, 1.029199559676785, 1.025433128034535, 1.032348017322074, 1.02956054640426, 1.021837572421634, 1.032589272239833, 1.02883621251074, 1.03572663446968, 1.01127052395519),
    @(1.02, 1.015757026439381, 1.029278809869973, 1.025656019899017, 1.032609451984648, 1.029561167926103, 1.021996893281358, 1.032626024211051, 1.029015921795615, 1.035945088389231, 1.011325256793188),
    @(1.02, 1.0158395106796, 1.029305822478734, 1.025732037653504, 1.03269861845784, 1.029561343914426, 1.02205121796653, 1.032638530682611, 1.029077203144091, 1.036019588056826, 1.011343915587419),
    @(1.02, 1.015470664081779, 1.029184977633886, 1.025392137006533, 1.032299939776394, 1.029560415031944, 1.021808266732241, 1.032582500137954, 1.028803158772577, 1.035686457458765, 1.011260454511378),
    @(1.02, 1.014272249036662, 1.028791441332969, 1.024288292544409, 1.031005459654759, 1.029554909994161, 1.021018435687276, 1.032398615060048, 1.027912569320712, 1.034604285185606, 1.010988859545931),
    @(1.02, 1.013519323194119, 1.028543518342299, 1.023595186264539, 1.030192842195708, 1.029549549582982, 1.020521858121802, 1.03228168691861, 1.027352892817892, 1.033924544930364, 1.010817902416904),
    @(1.02, 1.013918410019024, 1.028674992704325, 1.023962527342228, 1.030623505429295, 1.029552570889779, 1.020785101339317, 1.032343796520244, 1.027649562091004, 1.034284825157511, 1.010908548417934),
    @(1.02, 1.015490765290502, 1.029191566806252, 1.025410658770727, 1.03232166351239, 1.029560475047558, 1.021821508704644, 1.032585560613934, 1.028818094240572, 1.035704611492333, 1.011265004526705),
    @(1.02, 1.017317874621407, 1.029788797171083, 1.027095156568681, 1.034297833428938, 1.029561340424494, 1.023024291679433, 1.032860338570589, 1.030175296368112, 1.037355091090547, 1.011677799327959)
)

$columns = @("B","C","D","E","F","I","J","K","L","M","N")
$startRow = 2

for ($i = 0; $i -lt $vmPu.Count; $i++) {
    $row = $startRow + $i
    $values = $vmPu[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $ws.Range($columns[$j] + $row).Value = [double]$values[$j]
    }
}